# Insert a new weekly record for "Vega Modelo de Temuco" / "Frutilla" at row
# 140 (shifting the existing rows 140-165 down to 141-166), then fill in the
# new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 140..165 down to 141..166, leaving a blank row 140 behind.
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new data point.
$ws.Range("A140").Value = 10
$ws.Range("B140").Value = "Vega Modelo de Temuco"
$ws.Range("C140").Value = "La Araucanía"
$ws.Range("D140").Value = 44505
$ws.Range("E140").Value = 9
$ws.Range("F140").Value = "Fruta"
$ws.Range("G140").Value = 100101
$ws.Range("H140").Value = "Berries"
$ws.Range("I140").Value = 100112025
$ws.Range("J140").Value = "Frutilla"
$ws.Range("K140").Value = "Sin especificar"
$ws.Range("L140").Value = "Primera"
$ws.Range("M140").Value = 180
$ws.Range("N140").Value = 7000
$ws.Range("O140").Value = 7000
$ws.Range("P140").Value = 7000
$ws.Range("Q140").Value = "$/bandeja 7 kilos"
$ws.Range("R140").Value = "Provincia de Melipilla"
$ws.Range("S140").Value = 1000
$ws.Range("T140").Value = 7
